# Restored from revision of admin on 09/14/2020 07:08:42 AM.TEST
# Author: admin. Type: SAVE.
#
# Change: cell C10 on the "Rules" sheet (the "From" value of rule "R30")
# changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
